# Fruta / hortaliza, semanal
# Insert 3 new weekly records at rows 85-87 (shifting the existing
# rows 85-107 down to 88-110), and populate the new rows with data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows before the current row 85; this pushes the
# existing rows 85:107 down to 88:110 and keeps all of their data intact.
$ws.Rows("85:87").Insert()

# Columns that are constant for every record of this sheet/subset.
$mercadoId = 7
$mercado = "Terminal Hortofrutícola Agro Chillán"
$region = "Ñuble"
$codreg = 16
$categoriaId = 100112037
$categoria = "Cebollín"
$variedad = "Sin especificar"
$clasificacion = "Hortaliza"

# New row 85
$r = 85
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44985
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $categoriaId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = $variedad
$ws.Cells.Item($r, 9).Value = "Primera"
$ws.Cells.Item($r, 10).Value = 80
$ws.Cells.Item($r, 11).Value = 6000
$ws.Cells.Item($r, 12).Value = 6000
$ws.Cells.Item($r, 13).Value = 6000
$ws.Cells.Item($r, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item($r, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item($r, 16).Value = 167
$ws.Cells.Item($r, 17).Value = 36
$ws.Cells.Item($r, 18).Value = $clasificacion

# New row 86
$r = 86
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44985
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $categoriaId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = $variedad
$ws.Cells.Item($r, 9).Value = "Primera"
$ws.Cells.Item($r, 10).Value = 150
$ws.Cells.Item($r, 11).Value = 800
$ws.Cells.Item($r, 12).Value = 800
$ws.Cells.Item($r, 13).Value = 800
$ws.Cells.Item($r, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item($r, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item($r, 16).Value = 133
$ws.Cells.Item($r, 17).Value = 6
$ws.Cells.Item($r, 18).Value = $clasificacion

# New row 87
$r = 87
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44985
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $categoriaId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = $variedad
$ws.Cells.Item($r, 9).Value = "Segunda"
$ws.Cells.Item($r, 10).Value = 150
$ws.Cells.Item($r, 11).Value = 600
$ws.Cells.Item($r, 12).Value = 600
$ws.Cells.Item($r, 13).Value = 600
$ws.Cells.Item($r, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item($r, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item($r, 16).Value = 100
$ws.Cells.Item($r, 17).Value = 6
$ws.Cells.Item($r, 18).Value = $clasificacion
